# Apply the "slow dismantling expectations" test-data edit across several
# sheets in the AMIRIS scenario workbook: shift the simulation window back
# several years, rescale a handful of commodity/fuel prices, and replace
# the placeholder/unit-test power-plant rows with the real plant rows
# (dropping the now-unused dummy rows in the process).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "times": pull the simulation StartTime/StopTime back by 4 years.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("times")
$ws.Range("B2").Value = 43830.99861111111   # StartTime: 2023-12-31 -> 2019-12-31
$ws.Range("B3").Value = 44195.99861111111   # StopTime:  2024-12-30 -> 2020-12-30

# ---------------------------------------------------------------------
# Sheet "scenario_data_emlab": rebase the reference year + a few prices.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("scenario_data_emlab")
$ws.Range("B1").Value = 2020
$ws.Range("B2").Value = 20.4     # Co2Prices
$ws.Range("B5").Value = 10.8     # FuelPrice_HARD_COAL
$ws.Range("B6").Value = 20.16    # FuelPrice_NATURAL_GAS
$ws.Range("B7").Value = 46.44    # FuelPrice_OIL

# ---------------------------------------------------------------------
# Sheet "conventionals": drop the 5 extra NATURAL_GAS placeholder rows
# (old rows 8-12) and overwrite the remaining 6 data rows with the real
# plant data.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("conventionals")
$ws.Range("A8:A12").EntireRow.Delete()

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 19920300022
$ws.Range("C2").Value = "NATURAL_GAS"
$ws.Range("D2").Value = 4.2
$ws.Range("E2").Value = 0.61
$ws.Range("F2").Value = 31358.329
$ws.Range("G2").Value = 31358.329

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 19892800024
$ws.Range("C3").Value = "HARD_COAL"
$ws.Range("D3").Value = 3.5
$ws.Range("E3").Value = 0.33
$ws.Range("F3").Value = 24845.77
$ws.Range("G3").Value = 24845.77

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 19843000025
$ws.Range("C4").Value = "OIL"
$ws.Range("D4").Value = 6
$ws.Range("E4").Value = 0.35
$ws.Range("F4").Value = 3652.9
$ws.Range("G4").Value = 3652.9

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 19822900027
$ws.Range("C5").Value = "LIGNITE"
$ws.Range("D5").Value = 3.5
$ws.Range("E5").Value = 0.33
$ws.Range("F5").Value = 20779.02
$ws.Range("G5").Value = 20779.02

$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 19851400028
$ws.Range("C6").Value = "NUCLEAR"
$ws.Range("D6").Value = 3.5
$ws.Range("E6").Value = 0.33
$ws.Range("F6").Value = 8599
$ws.Range("G6").Value = 8599

$ws.Range("A7").Value = 5
$ws.Range("B7").Value = 19921700029
$ws.Range("C7").Value = "NATURAL_GAS"
$ws.Range("D7").Value = 4.5
$ws.Range("E7").Value = 0.43
$ws.Range("F7").Value = 8194.3025
$ws.Range("G7").Value = 8194.3025

# ---------------------------------------------------------------------
# Sheet "renewables": drop the 2 dummy placeholder rows (old rows 2-3)
# and the trailing row 6, shifting the real plant rows up; insert the
# new RunOfRiver row.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("renewables")
$ws.Range("A6").EntireRow.Delete()

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 20062400023
$ws.Range("C2").Value = 47547.50848700004
$ws.Range("D2").Value = 1.35
$ws.Range("E2").Value = "WindOn"

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 19641200026
$ws.Range("C3").Value = 8858.749999999998
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = "RunOfRiver"

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 20102100030
$ws.Range("C4").Value = 53555.51607579708
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = "OtherPV"

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 20142300031
$ws.Range("C5").Value = 10271.8
$ws.Range("D5").Value = 2.7
$ws.Range("E5").Value = "WindOff"

# ---------------------------------------------------------------------
# Sheet "biogas": drop the dummy placeholder row (old row 2) by removing
# the trailing row 3 and overwriting row 2 with the real plant data.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("biogas")
$ws.Range("A3").EntireRow.Delete()

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 20000100021
$ws.Range("C2").Value = 4644.4034
$ws.Range("D2").Value = 1.9
$ws.Range("E2").Value = "Biogas"
